$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.841.15'
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('D3').Value = '3.452.00'
$ws.Range('E3').Value = '  -1.18%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '159.52'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.89%  '
$ws.Range('D8').Value = '3.448.73'
$ws.Range('E8').Value = '  -1.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.577'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.94%  '
$ws.Range('E10').Value = '  -1.35%  '
$ws.Range('E12').Value = '  -1.39%  '
$ws.Range('D13').Value = '4.045.41'
$ws.Range('E13').Value = '  -1.12%  '
$ws.Range('E14').Value = '  -0.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.68'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000177'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -8.61%  '
$ws.Range('D17').Value = '64.888.02'
$ws.Range('E17').Value = '  -0.68%  '
$ws.Range('D18').Value = '3.429.79'
$ws.Range('E18').Value = '  -1.52%  '
$ws.Range('E19').Value = '  -3.69%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.77'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '378.33'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.85%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.97'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '72.22'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.88%  '
$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.535'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000122'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.91'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.66%  '
$ws.Range('E28').Value = '  -0.34%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.46'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.00%  '
$ws.Range('E31').Value = '  -2.40%  '
$ws.Range('E32').Value = '  -2.19%  '
$ws.Range('E33').Value = '  -2.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.85%  '
$ws.Range('E35').Value = '  -3.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '160.84'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.88'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.68%  '
$ws.Range('D38').Value = '2.898.83'
$ws.Range('E38').Value = '  -4.42%  '
$ws.Range('E39').Value = '  -3.84%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.24'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.57'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.64%  '
$ws.Range('E42').Value = '  -1.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.99'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.788'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.79%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '26.12'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0312'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.38'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +8.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '320.73'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.18%  '
$ws.Range('E49').Value = '  -2.76%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.48'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.24%  '
$ws.Range('E51').Value = '  -4.50%  '
